# Update the "Chart" sheet (GSC Coverage export) with the latest day-by-day
# crawl/index numbers pulled from the new export: the previously-missing
# 2025-10-24 row now has real data, the whole series shifts down one day,
# and three new trailing days (2026-01-18 .. 2026-01-20) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Column A holds ISO date strings ("yyyy-MM-dd"). Force the range to Text
# format first so Excel doesn't "helpfully" reinterpret them as date
# serials when we assign the string values below.
$ws.Range("A2:A90").NumberFormat = "@"

# r, Date, "Not indexed", "Indexed", "Impressions"
$rows = @(
    @(2, "2025-10-24", "__BLANK__", "__BLANK__", 11.0),
    @(3, "2025-10-25", 148.0, 136.0, 11.0),
    @(4, "2025-10-26", 132.0, 153.0, 8.0),
    @(5, "2025-10-27", 132.0, 153.0, 11.0),
    @(6, "2025-10-28", 132.0, 153.0, 20.0),
    @(7, "2025-10-29", 113.0, 181.0, 9.0),
    @(8, "2025-10-30", 113.0, 181.0, 13.0),
    @(9, "2025-10-31", 113.0, 181.0, 15.0),
    @(10, "2025-11-01", 113.0, 181.0, 21.0),
    @(11, "2025-11-02", 100.0, 203.0, 26.0),
    @(12, "2025-11-03", 100.0, 203.0, 15.0),
    @(13, "2025-11-04", 100.0, 203.0, 25.0),
    @(14, "2025-11-05", 102.0, 205.0, 31.0),
    @(15, "2025-11-06", 102.0, 205.0, 34.0),
    @(16, "2025-11-07", 102.0, 205.0, 23.0),
    @(17, "2025-11-08", 102.0, 205.0, 17.0),
    @(18, "2025-11-09", 106.0, 205.0, 19.0),
    @(19, "2025-11-10", 106.0, 205.0, 18.0),
    @(20, "2025-11-11", 106.0, 205.0, 26.0),
    @(21, "2025-11-12", 118.0, 202.0, 19.0),
    @(22, "2025-11-13", 118.0, 202.0, 13.0),
    @(23, "2025-11-14", 118.0, 202.0, 26.0),
    @(24, "2025-11-15", 118.0, 202.0, 18.0),
    @(25, "2025-11-16", 122.0, 208.0, 32.0),
    @(26, "2025-11-17", 122.0, 208.0, 33.0),
    @(27, "2025-11-18", 122.0, 208.0, 31.0),
    @(28, "2025-11-19", 127.0, 213.0, 24.0),
    @(29, "2025-11-20", 127.0, 213.0, 26.0),
    @(30, "2025-11-21", 127.0, 213.0, 38.0),
    @(31, "2025-11-22", 127.0, 213.0, 24.0),
    @(32, "2025-11-23", 180.0, 225.0, 33.0),
    @(33, "2025-11-24", 180.0, 225.0, 40.0),
    @(34, "2025-11-25", 180.0, 225.0, 46.0),
    @(35, "2025-11-26", 180.0, 225.0, 42.0),
    @(36, "2025-11-27", 180.0, 225.0, 47.0),
    @(37, "2025-11-28", 180.0, 225.0, 42.0),
    @(38, "2025-11-29", 180.0, 225.0, 24.0),
    @(39, "2025-11-30", 180.0, 225.0, 43.0),
    @(40, "2025-12-01", 180.0, 225.0, 53.0),
    @(41, "2025-12-02", 180.0, 225.0, 34.0),
    @(42, "2025-12-03", 180.0, 225.0, 35.0),
    @(43, "2025-12-04", 180.0, 225.0, 32.0),
    @(44, "2025-12-05", 180.0, 225.0, 29.0),
    @(45, "2025-12-06", 180.0, 225.0, 29.0),
    @(46, "2025-12-07", 180.0, 225.0, 36.0),
    @(47, "2025-12-08", 180.0, 225.0, 89.0),
    @(48, "2025-12-09", 180.0, 225.0, 59.0),
    @(49, "2025-12-10", 180.0, 225.0, 52.0),
    @(50, "2025-12-11", 180.0, 225.0, 62.0),
    @(51, "2025-12-12", 180.0, 225.0, 79.0),
    @(52, "2025-12-13", 180.0, 225.0, 53.0),
    @(53, "2025-12-14", 180.0, 225.0, 51.0),
    @(54, "2025-12-15", 180.0, 225.0, 55.0),
    @(55, "2025-12-16", 200.0, 231.0, 70.0),
    @(56, "2025-12-17", 200.0, 231.0, 54.0),
    @(57, "2025-12-18", 200.0, 231.0, 83.0),
    @(58, "2025-12-19", 200.0, 231.0, 63.0),
    @(59, "2025-12-20", 200.0, 231.0, 57.0),
    @(60, "2025-12-21", 200.0, 231.0, 54.0),
    @(61, "2025-12-22", 200.0, 231.0, 79.0),
    @(62, "2025-12-23", 200.0, 231.0, 54.0),
    @(63, "2025-12-24", 244.0, 227.0, 51.0),
    @(64, "2025-12-25", 244.0, 227.0, 50.0),
    @(65, "2025-12-26", 244.0, 227.0, 62.0),
    @(66, "2025-12-27", 244.0, 227.0, 59.0),
    @(67, "2025-12-28", 244.0, 227.0, 42.0),
    @(68, "2025-12-29", 244.0, 227.0, 40.0),
    @(69, "2025-12-30", 244.0, 227.0, 36.0),
    @(70, "2025-12-31", 244.0, 227.0, 72.0),
    @(71, "2026-01-01", 244.0, 227.0, 39.0),
    @(72, "2026-01-02", 244.0, 227.0, 46.0),
    @(73, "2026-01-03", 244.0, 227.0, 84.0),
    @(74, "2026-01-04", 246.0, 228.0, 90.0),
    @(75, "2026-01-05", 246.0, 228.0, 73.0),
    @(76, "2026-01-06", 246.0, 228.0, 79.0),
    @(77, "2026-01-07", 266.0, 233.0, 66.0),
    @(78, "2026-01-08", 266.0, 233.0, 58.0),
    @(79, "2026-01-09", 266.0, 233.0, 169.0),
    @(80, "2026-01-10", 266.0, 233.0, 183.0),
    @(81, "2026-01-11", 271.0, 232.0, 83.0),
    @(82, "2026-01-12", 271.0, 232.0, 92.0),
    @(83, "2026-01-13", 271.0, 232.0, 183.0),
    @(84, "2026-01-14", 281.0, 230.0, 136.0),
    @(85, "2026-01-15", 281.0, 230.0, 182.0),
    @(86, "2026-01-16", 281.0, 230.0, 140.0),
    @(87, "2026-01-17", 281.0, 230.0, 80.0),
    @(88, "2026-01-18", 292.0, 231.0, 50.0),
    @(89, "2026-01-19", 292.0, 231.0, 68.0),
    @(90, "2026-01-20", 292.0, 231.0, "__BLANK__")
)

foreach ($row in $rows) {
    $r = $row[0]
    $dateStr = $row[1]
    $notIndexed = $row[2]
    $indexed = $row[3]
    $impressions = $row[4]

    $ws.Cells.Item($r, 1).Value = $dateStr

    if ("$notIndexed" -ne "__BLANK__") {
        $ws.Cells.Item($r, 2).Value = $notIndexed
    }
    if ("$indexed" -ne "__BLANK__") {
        $ws.Cells.Item($r, 3).Value = $indexed
    }
    if ("$impressions" -ne "__BLANK__") {
        $ws.Cells.Item($r, 4).Value = $impressions
    }
}

# "Critical issues" sheet: refreshed Pages counts for a few reasons.
$wsCritical = $wb.Worksheets.Item("Critical issues")
$wsCritical.Cells.Item(2, 4).Value = 95.0
$wsCritical.Cells.Item(4, 4).Value = 59.0
$wsCritical.Cells.Item(6, 4).Value = 12.0
$wsCritical.Cells.Item(10, 4).Value = 17.0
